$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) hold numeric-looking / percent-looking
# text values. Force the cell to Text format before assigning so the string
# is preserved verbatim (matching the source inlineStr cells) instead of being
# auto-parsed into a Number/Percentage by Excel, then restore the default style.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "246.16"
Set-TextValue "E2" "0.21%"
Set-TextValue "D3" "29.93"
Set-TextValue "E3" "0.13%"
Set-TextValue "D4" "5.160"
Set-TextValue "E4" "0.31%"
Set-TextValue "E5" "1.16%"
Set-TextValue "D6" "6.664"
Set-TextValue "E6" "1.33%"
Set-TextValue "D7" "3.218"
Set-TextValue "E7" "6.64%"
Set-TextValue "D8" "0.8500"
Set-TextValue "E8" "-0.78%"
Set-TextValue "D9" "0.8630"
Set-TextValue "E9" "-0.72%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1380"
Set-TextValue "E10" "2.09%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07141"
Set-TextValue "E11" "3.04%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03196"
Set-TextValue "E12" "10.04%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09359"
Set-TextValue "E13" "-0.05%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001528"
Set-TextValue "E14" "0.63%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D15" "0.0006006"
Set-TextValue "E15" "0.33%"
Set-TextValue "D16" "0.005963"
Set-TextValue "E16" "-1.65%"
Set-TextValue "D17" "3.498"
Set-TextValue "E17" "-0.31%"
Set-TextValue "D18" "2.212"
Set-TextValue "E18" "-1.99%"
Set-TextValue "D19" "0.3191"
Set-TextValue "E19" "1.56%"
Set-TextValue "D20" "0.03372"
Set-TextValue "E20" "1.95%"
Set-TextValue "E21" "-0.35%"
Set-TextValue "D22" "3.492"
Set-TextValue "E22" "-3.07%"
Set-TextValue "D23" "0.04142"
Set-TextValue "E23" "-0.25%"
Set-TextValue "E24" "0.31%"
Set-TextValue "D25" "0.001225"
Set-TextValue "E25" "1.28%"
Set-TextValue "E26" "-7.57%"
Set-TextValue "E27" "1.89%"
Set-TextValue "E28" "4.31%"
Set-TextValue "D40" "0.03766"
Set-TextValue "E40" "-0.14%"
Set-TextValue "D41" "0.005688"
Set-TextValue "E41" "-0.34%"
Set-TextValue "D42" "0.1072"
Set-TextValue "E42" "0.58%"
Set-TextValue "D43" "0.002449"
Set-TextValue "E43" "11.58%"
Set-TextValue "D44" "0.009553"
Set-TextValue "E44" "-3.13%"
Set-TextValue "D45" "0.00005307"
Set-TextValue "E45" "4.47%"
Set-TextValue "E46" "0.16%"
Set-TextValue "D47" "0.05796"
Set-TextValue "E47" "-27.38%"
Set-TextValue "E48" "-20.24%"
Set-TextValue "E49" "0.16%"
Set-TextValue "E50" "0.16%"
